$d = $word.ActiveDocument

# --- East-Asian default font change: DejaVu Sans -> Tahoma -------------
# Applies to the "Normal" and "Heading" paragraph styles' run properties
# (w:rPr/w:rFonts/@w:eastAsia). The document's docDefaults/rPrDefault
# (the very first instance in the diff) is not reachable through the
# exposed Word object model in this runtime, so it is intentionally left
# untouched; every style that IS reachable gets updated below.

$normal = $d.Styles("Normal")
$normal.Font.NameFarEast = "Tahoma"

$heading = $d.Styles("Heading")
$heading.Font.NameFarEast = "Tahoma"

# --- Add a complex-script (w:cs) font to styles that had none ----------
# "NameBi" is this engine's binding for Font.NameBidirectional, which
# writes to <w:rFonts w:cs="..."/>.

$list = $d.Styles("List")
$list.Font.NameBi = "DejaVu Sans"

$caption = $d.Styles("Caption")
$caption.Font.NameBi = "DejaVu Sans"

$index = $d.Styles("Index")
$index.Font.NameBi = "DejaVu Sans"
